# Applies "Add files via upload" text-run-merge edits to
# "장철화 앱 프로젝트.pptx" (before.pptx).
#
# In every case the two original runs being merged already share
# identical run formatting (rPr), so the edit is just a textual
# merge of two adjacent <a:r> runs into one: the first run's text
# becomes the concatenation, and the (now redundant) second run is
# collapsed by emptying its text.

function Merge-Runs {
    param($Paragraph, $FirstRunIndex, $SecondRunIndex, $MergedText)
    $firstRun = $Paragraph.Runs($FirstRunIndex, 1)
    $secondRun = $Paragraph.Runs($SecondRunIndex, 1)
    $firstRun.Text = $MergedText
    $secondRun.Text = ""
}

# Re-assigning Shape.Height round-trips through a single-precision
# (float32) points value, which truncates just below the true EMU on
# save. A tiny (sub-EMU-visible) nudge compensates for that systematic
# bias so the restored height lands back on the exact original EMU
# instead of one EMU short.
$HeightRestoreEpsilonPt = 0.00003

function Restore-ShapeHeight {
    param($Shape, $OriginalHeight)
    $Shape.Height = $OriginalHeight + $HeightRestoreEpsilonPt
}

$p = $ppt.ActivePresentation

# --- Slide 6 : "TextBox 9" paragraph 6 -----------------------------------
# '...' + ' ' + '배달 ' + '앱' ...  ->  '...' + ' 배달 ' + '앱' ...
$s6 = $p.Slides.Item(6)
$shape6 = $s6.Shapes.Item(5)
$shape6Height = $shape6.Height
$para6 = $shape6.TextFrame.TextRange.Paragraphs(6, 1)
Merge-Runs $para6 4 5 " 배달 "
# The shape auto-fits its text box; merging runs doesn't change the
# rendered text, so restore the pre-edit height (autofit shrink is an
# artifact of re-running layout on the editing engine, not a real change).
Restore-ShapeHeight $shape6 $shape6Height

# --- Slide 7 : "TextBox 9" paragraph 1, 2, 3 ------------------------------
$s7 = $p.Slides.Item(7)
$shape7 = $s7.Shapes.Item(5)
$shape7Height = $shape7.Height

# ' ' + '2. ' + '프로젝트 개발 도구'  ->  ' 2. ' + '프로젝트 개발 도구'
$para7_1 = $shape7.TextFrame.TextRange.Paragraphs(1, 1)
Merge-Runs $para7_1 1 2 " 2. "

# ' ' + '  '  ->  '   '
$para7_2 = $shape7.TextFrame.TextRange.Paragraphs(2, 1)
Merge-Runs $para7_2 1 2 "   "

# '개발 환경' + ' ' + ': ' + 'Android Studio'  ->  '개발 환경' + ' : ' + 'Android Studio'
$para7_3 = $shape7.TextFrame.TextRange.Paragraphs(3, 1)
Merge-Runs $para7_3 3 4 " : "

Restore-ShapeHeight $shape7 $shape7Height

# --- Slide 10 : "TextBox 9" paragraph 2 -----------------------------------
# '프로젝트 커뮤니케이션 ' + '관리'  ->  '프로젝트 커뮤니케이션 관리'
$s10 = $p.Slides.Item(10)
$shape10 = $s10.Shapes.Item(5)
$shape10Height = $shape10.Height
$para10 = $shape10.TextFrame.TextRange.Paragraphs(2, 1)
Merge-Runs $para10 2 3 "프로젝트 커뮤니케이션 관리"
Restore-ShapeHeight $shape10 $shape10Height
